$d = $word.ActiveDocument

# The empty paragraph right after "...Retirar "_rc" do final" (the first of
# three blank placeholder paragraphs that precede "Arquivos:") is turned
# into a new bulleted section heading "Compilação do programa para .exe:",
# followed by a blank spacer paragraph and two new body paragraphs
# explaining the .exe compilation step.
$anchor = $d.Content
[void]$anchor.Find.Execute("do final", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hitPara = $anchor.Paragraphs(1)
$target = $hitPara.Next()
$r = $target.Range

$xml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:cstheme="minorHAnsi"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:cstheme="minorHAnsi"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>Compilação do programa para .exe</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:cstheme="minorHAnsi"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:cstheme="minorHAnsi"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:ind w:firstLine="708"/>
              <w:jc w:val="both"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Para realizar a compilação dos programas é necessário rodar um comando do </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>pyinstaller</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> apropriado para cada um deles. No código principal se tem o comando apropriado ou como se inspirar para criar o comando apropriado. </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:ind w:firstLine="708"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:cstheme="minorHAnsi"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
            </w:pPr>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>Obs</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>: algumas vezes é necessário adicionar uma pasta ou um arquivo para que se funcione adequadamente.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

[void]$r.InsertXML($xml)
